# Cholesky unit now passes validataion.
#
# Reproduces the workbook/sheet view tweaks and the formula fix captured
# in the target diff:
#   1. Reposition the workbook window (bookViews/workbookView).
#   2. Scroll Sheet1's frozen pane so its top-left visible cell is E17.
#   3. Update the active selection in the bottom-right pane to W24.
#   4. Fix cell E19's formula (MAX(2:2) -> MAX(2:2)+1); Excel recalculates
#      the shared "+1" formula chain across F19:DM19 automatically.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Sheet1")
$win = $excel.ActiveWindow

# --- 1. Workbook window position ---------------------------------------
$win.Left = -38400
$win.Top  = 600

# --- 2. Sheet1: scroll the frozen pane to E17 ---------------------------
$ws.Activate() | Out-Null

$win.ScrollColumn = 5   # column E
$win.ScrollRow    = 17  # row 17

# --- 3. Sheet1: set the bottom-right pane's active selection ------------
$ws.Range("W24").Select() | Out-Null

# --- 4. Fix the formula in E19 (shared dependents F19:DM19 recalc) ------
$ws.Range("E19").Formula = "=MAX(2:2)+1"
